$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hours spent on tasks (columns D/E/F) for several rows
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 2

$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 2

$ws.Range("F14").Value = 2

$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 2

$ws.Range("E16").Value = 2

$ws.Range("D17").Value = 3

$ws.Range("D18").Value = 3

$ws.Range("D19").Value = 3

$ws.Range("D20").Value = 3

$excel.CalculateFull()

$wb.Save()
